$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:D values for sheet rows 2-67 (rows 2-57 updated in place, rows 58-67 newly added)
$bcdValues = @(
    @(0.5520055958676888,7.456712159883656,0.630815427435628),
    @(1.529003317215714,6.993709231954493,0.5932092276256263),
    @(2.806415998886871,4.550745007401177,0.3446298390482985),
    @(3.088822840973728,8.829843762265867,0.5025077914174758),
    @(5.354586240738302,4.383563488021912,0.4871749539679492),
    @(6.447604160660143,2.674706966612216,0.3347607981322681),
    @(6.77780262091801,2.562529901237825,0.2860092982698995),
    @(6.850577864769429,4.8755504381315,0.8272331001450209),
    @(7.218474110226925,6.335385122369471,0.4680818590937341),
    @(7.660797315583311,3.414679107069724,0.4133327348899845),
    @(10.64360685248097,7.395326987564623,0.4824865444029217),
    @(12.42132274825905,4.436952232130229,0.3834920061373053),
    @(13.46863987586207,4.312380270511887,0.6009280333151326),
    @(15.67551034550826,2.876710794586397,0.4521522995427405),
    @(16.47357041818529,1.321099829281096,0.1795779920549424),
    @(18.95607699594084,7.930374260151479,0.6337808690348343),
    @(20.89317459886447,2.565235453999529,0.2515505572437173),
    @(21.09840681968731,4.25655825314167,0.3010546989728326),
    @(25.30680571927003,4.368743629561925,0.6288879214942688),
    @(26.16524394627741,3.585656057004245,0.5916584004188435),
    @(26.3811607436705,3.921779212430995,0.4663696702334351),
    @(33.10898143613653,5.22309853976534,0.2367969411308504),
    @(34.4607725385481,3.578481134717975,0.3448723589882984),
    @(36.24780611978512,3.676351811910774,0.4537934267754338),
    @(37.18140980391463,3.042433037719632,0.247805941979705),
    @(38.99508142126555,6.516189257080084,0.3010958393371924),
    @(39.54769663529568,6.781946532068721,0.6470481662673447),
    @(39.55121359874423,3.306607903235394,0.4293313467907782),
    @(39.8101834115626,6.753064993601121,0.7063965314361182),
    @(40.04983599306401,2.938168133706037,0.4541418681362107),
    @(40.34907288398516,3.339458623601474,0.3130155285277615),
    @(41.34062858174761,4.610220271458755,0.5311268357233145),
    @(44.22973943203641,5.884998036996139,0.8343037225370524),
    @(44.30724434226106,2.995979210580881,0.5691880264284979),
    @(47.21608749573338,2.62790713009155,0.4744401093592391),
    @(47.60416042484564,5.739313250170197,0.8424123423958878),
    @(49.34469191643922,7.664094613157691,0.7663935726074347),
    @(49.65908855196219,5.899644030180284,0.6699558081458459),
    @(50.7230050865786,5.803548391683157,0.5988114634904502),
    @(54.21178607093759,4.956154329027915,0.7146003377038491),
    @(55.74850404702884,3.528097592827379,0.4588315675156596),
    @(58.44957369345386,7.210811298644916,0.7581691955676273),
    @(60.14409442043063,7.840228315861816,0.7811369341510174),
    @(60.30724993710545,6.167867675051744,0.6356133484366381),
    @(63.66848422183852,7.080550686794093,0.2350020873408501),
    @(66.76608104416796,5.34972229333688,0.539874398876771),
    @(69.53986120974791,8.179125219367217,0.6431927438960214),
    @(69.81748527331543,3.705529363014854,0.1749910638023167),
    @(72.04859277444831,7.436864957084602,0.7300699733760003),
    @(76.12353685879624,4.506241233691332,0.7626794478801696),
    @(76.70978554460662,3.997235605461916,0.7496416425416778),
    @(77.05857018353711,6.290392578644001,0.5427019548473698),
    @(78.39984561514854,4.473980890713594,0.5943271270065704),
    @(78.47587246877349,2.587438704591231,0.8066445593957717),
    @(78.80357509137781,2.445890147298472,0.3427389864343864),
    @(79.12059466182882,5.937740844124311,0.6532655772524482),
    @(80.72791339274384,1.83553724674939,0.3253029172779324),
    @(84.33025024973355,2.155515905017443,0.504307006067259),
    @(88.19724037068751,7.927903183158886,0.6069506079480076),
    @(88.89511370022807,1.73724213554368,0.3058145734594823),
    @(89.763102968146,3.460911710339154,0.5176506839751086),
    @(92.16975532985956,4.385905742323472,0.601422052032637),
    @(96.7062391798861,3.717756753320336,0.7578386739169989),
    @(97.41297116665748,5.496958809094511,0.6427022769806511),
    @(98.01449985644851,1.622490803369615,0.2229765514957482),
    @(99.20463479035409,2.422020837435986,0.5708743424287757)
)

for ($i = 0; $i -lt $bcdValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $bcdValues[$i][0]
    $ws.Cells.Item($row, 3).Value = $bcdValues[$i][1]
    $ws.Cells.Item($row, 4).Value = $bcdValues[$i][2]
}

# New A-column values for the newly added rows 58-67 (continuing the 0-based index sequence),
# copying the existing A-column cell format (bordered, centered, bold) onto each new cell.
$aValues = @(56,57,58,59,60,61,62,63,64,65)
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = 58 + $i
    $ws.Range("A2").Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
}

$ws.Range("A1").Select()
Write-Host "done"
